$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (G1) to the
# new header cell (H1), so the new "Save" header matches the styling used by
# all the other header cells (bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and the corresponding numeric value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
